$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(563).Insert()
$ws.Range("A563").NumberFormat = "@"
$ws.Range("A563").Value = "2026/01/05"
$ws.Range("A563").Style = "Normal"
$ws.Range("B563").Value = "月"
$ws.Range("C563").Value = 13
$ws.Range("D563").Value = 144
